$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (col 8); existing H..P shift right to I..Q
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column
$ws.Cells.Item(1, 8).Value = "CO2/(CO+CO2)"

# Give the new column's cells the same border/number-format "style" as the
# other plain data columns (E:G) -- this reuses the existing style index
# rather than minting a new one.
$ws.Range("H1:H31").Borders.LineStyle = 1

# Row 2's formula is entered standalone first ...
$ws.Range("H2").Formula = "=F2/(E2+F2)"
# ... then rows 3-31 are filled as a separate operation so Excel groups them
# into one shared-formula block (H3:H31) distinct from H2.
$ws.Range("H3:H31").Formula = "=F3/(E3+F3)"

# Match the author's final selection
$ws.Range("H3").Select()
